$d = $word.ActiveDocument

# Update the date line (first paragraph)
$d.Paragraphs(1).Range.Text = "2026-01-01 Thursday"

# Update each answer cell in the table, row-major order
$newValues = @("66-35=31","82-46=36","49-11=38","52+5=57","95-4=91","14+62=76","44+50=94","64+19=83","30+16=46","60+26=86","56+19=75","6+13=19","38-27=11","96-51=45","19+44=63","13+7=20","59+16=75","27+52=79","41-9=32","52-26=26","36+16=52","50-36=14","44+17=61","73-29=44","50-13=37","87-67=20","82-65=17","50-50=0","46+29=75","91-65=26","89-25=64","3+55=58","72+11=83","20+76=96","82+7=89","12+0=12","77-14=63","6-0=6","14+63=77","93-24=69","73-57=16","26+70=96","12-8=4","52+1=53","21+34=55","83+4=87","3+26=29","54-6=48","95-89=6","13+6=19","17+26=43","12+46=58","56-8=48","79-60=19","19-13=6","32+40=72","57-23=34","44-10=34","5+70=75","3+40=43","1+60=61","32-4=28","4+26=30","2+84=86","28+9=37","5+87=92","29+51=80","30-6=24","69-43=26","13+53=66","66-63=3","70+13=83","36-30=6","80-70=10","4+81=85","14-12=2","57-36=21","55-43=12","55-12=43","94-82=12","76-58=18","24-1=23","81+7=88","9+7=16","68-9=59","64-14=50","81-26=55","46+0=46","33-10=23","36-2=34","88-86=2","76-22=54","78-74=4","8+73=81","37-18=19","98-54=44","80-67=13","46+48=94","56+3=59","23+10=33")

$t = $d.Tables(1)
$rows = $t.Rows.Count
$cols = $t.Columns.Count

$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Output "done: $idx cells updated"
